# Add a new "Locator Type" column (D) that classifies each existing
# "Locator Value" (column C) as either an Xpath expression (starts with
# "//") or a CSS selector (everything else), to support running tests on
# the Android mobile browser.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 49

# Header cell - gets its own distinct font (Helvetica Neue, 11pt).
$headerCell = $ws.Cells.Item(1, 4)
$headerCell.Value = "Locator Type"
$headerCell.Font.Name = "Helvetica Neue"
$headerCell.Font.Size = 11

for ($r = 2; $r -le $lastRow; $r++) {
    $locatorValue = $ws.Cells.Item($r, 3).Text

    if ($locatorValue.StartsWith("//")) {
        $locatorType = "Xpath"
    } else {
        $locatorType = "CSS"
    }

    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = $locatorType

    # Rows below 30 previously had no D/E cell at all, so a brand-new cell
    # picks up the column's default style instead of the plain-data style
    # used throughout the rest of the sheet. Force the same font the rest
    # of the table's un-styled cells use (Helvetica 11) so the cell style
    # lines up with the existing D2:D30 cells.
    $cell.Font.Name = "Helvetica"
    $cell.Font.Size = 11
}

# Widen column C to fit the longer content and add an explicit width for
# the new column D.
$ws.Columns.Item(3).ColumnWidth = 85.67
$ws.Columns.Item(4).ColumnWidth = 14

# Best-effort cosmetic view updates (scroll position / selection) to
# mirror where the author was working when the column was added.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 20
$ws.Range("C46").Select() | Out-Null
